$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 6536282
$ws.Range("I41").Value = 7936828.5
$ws.Range("K41").Value = 7936828.5
$ws.Range("M41").Value = -7936388.5
$ws.Range("H43").Value = 1048.6666
$ws.Range("I43").Value = 733.3333
$ws.Range("J43").Value = 1364
$ws.Range("K43").Value = 733.3333
$ws.Range("L43").Value = 1364
$ws.Range("M43").Value = -664.3333
$ws.Range("N43").Value = -1502
$ws.Range("H64").Value = 7817769
$ws.Range("I64").Value = 20836666
$ws.Range("J64").Value = 6430.3
$ws.Range("K64").Value = 20836666
$ws.Range("L64").Value = 6430.3
$ws.Range("M64").Value = -20836418
$ws.Range("N64").Value = -6926.3
$ws.Range("H67").Value = 7817769
$ws.Range("I67").Value = 20836666
$ws.Range("J67").Value = 6430.3
$ws.Range("K67").Value = 20836666
$ws.Range("L67").Value = 6430.3
$ws.Range("M67").Value = -20835808
$ws.Range("N67").Value = -8146.3
$ws.Range("H74").Value = 4173.091
$ws.Range("I74").Value = 3766.6667
$ws.Range("K74").Value = 3766.6667
$ws.Range("M74").Value = -2830.6667
$ws.Range("H76").Value = 3971193
$ws.Range("I76").Value = 4447324
$ws.Range("J76").Value = 3433.3333
$ws.Range("K76").Value = 4447324
$ws.Range("L76").Value = 3433.3333
$ws.Range("M76").Value = -4447009
$ws.Range("N76").Value = -4063.3333
$ws.Range("H77").Value = 4173.091
$ws.Range("I77").Value = 3766.6667
$ws.Range("K77").Value = 18833.3335
$ws.Range("M77").Value = -14153.3335
$ws.Range("H79").Value = 3971193
$ws.Range("I79").Value = 4447324
$ws.Range("J79").Value = 3433.3333
$ws.Range("K79").Value = 4447324
$ws.Range("L79").Value = 3433.3333
$ws.Range("M79").Value = -4446232
$ws.Range("N79").Value = -5617.3333
$ws.Range("H129").Value = 1016.6774
$ws.Range("I129").Value = 265.66666
$ws.Range("J129").Value = 1097.1428
$ws.Range("K129").Value = 796.9999799999999
$ws.Range("L129").Value = 3291.4284
$ws.Range("M129").Value = 4203.00002
$ws.Range("N129").Value = -13291.4284
$ws.Range("H133").Value = 49485
$ws.Range("J133").Value = 49485
$ws.Range("L133").Value = 49485
$ws.Range("N133").Value = -59605
$ws.Range("H137").Value = 47621090
$ws.Range("I137").Value = 71429960
$ws.Range("J137").Value = 3343.7144
$ws.Range("K137").Value = 214289880
$ws.Range("L137").Value = 10031.1432
$ws.Range("M137").Value = -214287330
$ws.Range("N137").Value = -15131.1432
$ws.Range("H141").Value = 2027.5
$ws.Range("I141").Value = 1970.2941
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 5910.8823
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -730.8823000000002
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 36463.48
$ws.Range("I32").Value = 9615.473
$ws.Range("J32").Value = 105501.21
$ws.Range("K32").Value = 9615.473
$ws.Range("L32").Value = 105501.21
$ws.Range("M32").Value = -9328.473
$ws.Range("N32").Value = -106075.21
$ws.Range("H63").Value = 16278.214
$ws.Range("I63").Value = 19988.334
$ws.Range("J63").Value = 9600
$ws.Range("K63").Value = 19988.334
$ws.Range("L63").Value = 9600
$ws.Range("M63").Value = -19302.334
$ws.Range("N63").Value = -10972
$ws.Range("H66").Value = 16278.214
$ws.Range("I66").Value = 19988.334
$ws.Range("J66").Value = 9600
$ws.Range("K66").Value = 99941.67
$ws.Range("L66").Value = 48000
$ws.Range("M66").Value = -96509.67
$ws.Range("N66").Value = -54864
$ws.Range("H74").Value = 5326.75
$ws.Range("I74").Value = 1029.1818
$ws.Range("J74").Value = 14781.4
$ws.Range("K74").Value = 1029.1818
$ws.Range("L74").Value = 14781.4
$ws.Range("M74").Value = -155.1818000000001
$ws.Range("N74").Value = -16529.4
$ws.Range("H77").Value = 5326.75
$ws.Range("I77").Value = 1029.1818
$ws.Range("J77").Value = 14781.4
$ws.Range("K77").Value = 5145.909000000001
$ws.Range("L77").Value = 73907
$ws.Range("M77").Value = -777.9090000000006
$ws.Range("N77").Value = -82643
$ws.Range("H88").Value = 3550
$ws.Range("J88").Value = 2260
$ws.Range("L88").Value = 2260
$ws.Range("N88").Value = -3072
$ws.Range("H91").Value = 3550
$ws.Range("J91").Value = 2260
$ws.Range("L91").Value = 2260
$ws.Range("N91").Value = -5068
$ws.Range("H132").Value = 1947.5518
$ws.Range("I132").Value = 1819.7551
$ws.Range("J132").Value = 2643.3333
$ws.Range("K132").Value = 5459.2653
$ws.Range("L132").Value = 7929.999899999999
$ws.Range("M132").Value = -2929.2653
$ws.Range("N132").Value = -12989.9999
$ws.Range("H133").Value = 40713.445
$ws.Range("J133").Value = 40713.445
$ws.Range("L133").Value = 40713.445
$ws.Range("N133").Value = -45773.445
$ws.Range("H139").Value = 54857.5
$ws.Range("J139").Value = 54857.5
$ws.Range("L139").Value = 54857.5
$ws.Range("N139").Value = -65137.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 550
$ws.Range("I8").Value = 200
$ws.Range("J8").Value = 1250
$ws.Range("K8").Value = 200
$ws.Range("L8").Value = 1250
$ws.Range("M8").Value = -60
$ws.Range("N8").Value = -1530
$ws.Range("H59").Value = 59950
$ws.Range("J59").Value = 59950
$ws.Range("L59").Value = 59950
$ws.Range("N59").Value = -61644
$ws.Range("H105").Value = 2816.7827
$ws.Range("I105").Value = 2695.6562
$ws.Range("J105").Value = 3093.6428
$ws.Range("K105").Value = 2695.6562
$ws.Range("L105").Value = 3093.6428
$ws.Range("M105").Value = -948.6561999999999
$ws.Range("N105").Value = -6587.6428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 26317408
$ws.Range("I58").Value = 30304328
$ws.Range("J58").Value = 3728.4
$ws.Range("K58").Value = 30304328
$ws.Range("L58").Value = 3728.4
$ws.Range("M58").Value = -30304125
$ws.Range("N58").Value = -4134.4
$ws.Range("H62").Value = 16004.723
$ws.Range("I62").Value = 18234.643
$ws.Range("K62").Value = 18234.643
$ws.Range("M62").Value = -17610.643
$ws.Range("H65").Value = 16004.723
$ws.Range("I65").Value = 18234.643
$ws.Range("K65").Value = 91173.215
$ws.Range("M65").Value = -88053.215
$ws.Range("H136").Value = 26317408
$ws.Range("I136").Value = 30304328
$ws.Range("J136").Value = 3728.4
$ws.Range("K136").Value = 90912984
$ws.Range("L136").Value = 11185.2
$ws.Range("M136").Value = -90910434
$ws.Range("N136").Value = -16285.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 6937.6562
$ws.Range("I134").Value = 3471.0588
$ws.Range("J134").Value = 10866.467
$ws.Range("K134").Value = 10413.1764
$ws.Range("L134").Value = 32599.401
$ws.Range("M134").Value = -5343.1764
$ws.Range("N134").Value = -42739.401

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5280.5137
$ws.Range("I70").Value = 5262.1562
$ws.Range("J70").Value = 5398
$ws.Range("K70").Value = 5262.1562
$ws.Range("L70").Value = 5398
$ws.Range("M70").Value = -4992.1562
$ws.Range("N70").Value = -5938
$ws.Range("H73").Value = 5280.5137
$ws.Range("I73").Value = 5262.1562
$ws.Range("J73").Value = 5398
$ws.Range("K73").Value = 5262.1562
$ws.Range("L73").Value = 5398
$ws.Range("M73").Value = -4326.1562
$ws.Range("N73").Value = -7270
$ws.Range("H80").Value = 65219990
$ws.Range("I80").Value = 2675
$ws.Range("J80").Value = 187502460
$ws.Range("K80").Value = 2675
$ws.Range("L80").Value = 187502460
$ws.Range("M80").Value = -1677
$ws.Range("N80").Value = -187504456
$ws.Range("H83").Value = 65219990
$ws.Range("I83").Value = 2675
$ws.Range("J83").Value = 187502460
$ws.Range("K83").Value = 13375
$ws.Range("L83").Value = 937512300
$ws.Range("M83").Value = -8383
$ws.Range("N83").Value = -937522284
$ws.Range("H107").Value = 820.8823
$ws.Range("I107").Value = 1125.8889
$ws.Range("J107").Value = 477.75
$ws.Range("K107").Value = 1125.8889
$ws.Range("L107").Value = 477.75
$ws.Range("M107").Value = 794.1111000000001
$ws.Range("N107").Value = -4317.75
$ws.Range("H122").Value = 2734
$ws.Range("I122").Value = 2742.577
$ws.Range("J122").Value = 2678.25
$ws.Range("K122").Value = 8227.731
$ws.Range("L122").Value = 8034.75
$ws.Range("M122").Value = -5777.731
$ws.Range("N122").Value = -12934.75
$ws.Range("H132").Value = 4108.3335
$ws.Range("I132").Value = 4403.3887
$ws.Range("J132").Value = 2338
$ws.Range("K132").Value = 13210.1661
$ws.Range("L132").Value = 7014
$ws.Range("M132").Value = -10680.1661
$ws.Range("N132").Value = -12074
$ws.Range("H138").Value = 93000
$ws.Range("J138").Value = 93000
$ws.Range("L138").Value = 93000
$ws.Range("N138").Value = -103280
$ws.Range("H139").Value = 59993.75
$ws.Range("J139").Value = 59993.75
$ws.Range("L139").Value = 59993.75
$ws.Range("N139").Value = -70273.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8050.7144
$ws.Range("I22").Value = 1147.7778
$ws.Range("J22").Value = 20476
$ws.Range("K22").Value = 1147.7778
$ws.Range("L22").Value = 20476
$ws.Range("M22").Value = -852.7778000000001
$ws.Range("N22").Value = -21066
$ws.Range("H26").Value = 50000
$ws.Range("I26").Value = 50000
$ws.Range("K26").Value = 50000
$ws.Range("M26").Value = -49705
$ws.Range("H27").Value = 8050.7144
$ws.Range("I27").Value = 1147.7778
$ws.Range("J27").Value = 20476
$ws.Range("K27").Value = 1147.7778
$ws.Range("L27").Value = 20476
$ws.Range("M27").Value = -1040.7778
$ws.Range("N27").Value = -20690
$ws.Range("H34").Value = 36666.668
$ws.Range("I34").Value = 36666.668
$ws.Range("K34").Value = 36666.668
$ws.Range("M34").Value = -36494.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 9375
$ws.Range("J20").Value = 9375
$ws.Range("L20").Value = 9375
$ws.Range("N20").Value = -9855
